# Apply the changes described in the commit:
# "#2 Commited test scripts, test outputs and test case for 7.0"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (Test Case #7): now marked as "Fail" with no justification text,
# and the row height shrinks since the long justification text is gone.
$ws.Range("G8").Value = "Fail"
$ws.Range("H8").Value = "-"
$ws.Rows.Item(8).RowHeight = 55.5

# Remove the screenshot/picture that was attached to the sheet.
for ($i = $ws.Shapes.Count; $i -ge 1; $i--) {
    $ws.Shapes.Item($i).Delete() | Out-Null
}

# Update the current selection/view to reflect where the author left off.
$ws.Range("H9").Select() | Out-Null
